$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: "No." summary row ---
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("B10").Value = 21
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 28

# --- Row 11: "Marking" summary row ---
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# --- Row 12: "Total" summary row ---
$ws.Range("A12").Style = "mtitleStyle"
$ws.Range("B12").Value = 84
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "83/112"

# --- Row 15: remove the third (G/H) "Student Ans"/"Correct Ans" header block ---
$ws.Range("G15:H15").Clear()

# --- Remove the third (G/H) answer-option block for rows 16-21 ---
$ws.Range("G16:H21").Clear()

# --- Fill column A ("Student Ans") with the student's answer, colored by correctness ---
$ws.Range("A16").Value = "Option A"
$ws.Range("A16").Style = "correctStyle"
$ws.Range("A17").Value = "Option D"
$ws.Range("A17").Style = "correctStyle"
$ws.Range("A18").Value = "Option B"
$ws.Range("A18").Style = "correctStyle"
$ws.Range("A19").Value = "Option C"
$ws.Range("A19").Style = "correctStyle"
$ws.Range("A20").Value = "Option B"
$ws.Range("A20").Style = "correctStyle"
$ws.Range("A21").Value = "Option C"
$ws.Range("A21").Style = "correctStyle"
$ws.Range("A22").Value = "Option D"
$ws.Range("A22").Style = "correctStyle"
$ws.Range("A25").Value = "Option A"
$ws.Range("A25").Style = "correctStyle"
$ws.Range("A26").Value = "Option C"
$ws.Range("A26").Style = "correctStyle"
$ws.Range("A27").Value = "Option A"
$ws.Range("A27").Style = "correctStyle"
$ws.Range("A29").Value = "Option D"
$ws.Range("A29").Style = "correctStyle"
$ws.Range("A31").Value = "Option D"
$ws.Range("A31").Style = "correctStyle"
$ws.Range("A32").Value = "Option C"
$ws.Range("A32").Style = "correctStyle"
$ws.Range("A33").Value = "Option D"
$ws.Range("A33").Style = "correctStyle"
$ws.Range("A34").Value = "Option A"
$ws.Range("A34").Style = "incorrectStyle"
$ws.Range("A36").Value = "Option A"
$ws.Range("A36").Style = "correctStyle"
$ws.Range("A38").Value = "Option A"
$ws.Range("A38").Style = "correctStyle"
$ws.Range("A39").Value = "Option D"
$ws.Range("A39").Style = "correctStyle"
$ws.Range("A40").Value = "Option D"
$ws.Range("A40").Style = "correctStyle"

# --- Fill column D ("Student Ans", second block) for rows 16-18 ---
$ws.Range("D16").Value = "Option A"
$ws.Range("D16").Style = "correctStyle"
$ws.Range("D17").Value = "Option C"
$ws.Range("D17").Style = "correctStyle"
$ws.Range("D18").Value = "Option D"
$ws.Range("D18").Style = "correctStyle"

# --- Remove the second (D/E) block for rows 19-40 (no longer needed) ---
$ws.Range("D19:E40").Clear()

